$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report date-range header (row 2, merged A2:D2)
$ws.Range("A2").Value = "گزارش کار از تاریخ 13 شهریور تا تاریخ 19 شهریور"

# Row 4 - Saturday (13 شهریور)
$ws.Range("A4").Value = "13 شهریور"
$ws.Range("C4").Value = "درست کردن استایل سایت، تلاش برای وصل کردن API، دیدن ویدیو آموزشی"
$ws.Range("D4").Value = 0.16666666666666666

# Row 5 - Sunday (شهریور 14)
$ws.Range("A5").Value = "شهریور 14"
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = 0

# Row 6 - Monday (شهریور 15)
$ws.Range("A6").Value = "شهریور 15"
$ws.Range("C6").Value = "جلسه کارآموزی، وصل کردن API، برطرف کردن باگ های دیزاین"
$ws.Range("D6").Value = 0.1875

# Row 7 - Tuesday (شهریور 16)
$ws.Range("A7").Value = "شهریور 16"
$ws.Range("C7").Value = "کار کردن روی وبسایت خیرین و تصحیح فولدر بندی ها"
$ws.Range("D7").Value = 0.10416666666666667

# Row 8 - Wednesday (شهریور 17)
$ws.Range("A8").Value = "شهریور 17"
$ws.Range("C8").Value = "نهایی کردن استایل وبسایت"
$ws.Range("D8").Value = 0.041666666666666664

# Row 9 - Thursday (شهریور 18)
$ws.Range("A9").Value = "شهریور 18"
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = 0

# Row 10 - Friday (شهریور 19)
$ws.Range("A10").Value = "شهریور 19"
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = 0

# Total hours worked
$ws.Range("D35").Value = " 11:00:00"

# Update view: scroll back to top and select A13
$ws.Range("A13").Select()
